# Add a new registrant row (row 4) to the follower sheet, duplicating the
# values of row 3 (same record) but with sequence number 3, and move the
# active cell selection to J2 (matches the authored diff for sheet1.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight       = -4152
$xlLeft        = -4131
$xlPasteValues = -4163

# --- styles first --------------------------------------------------------
# Columns A..Q (except J) and S..AI mirror the "right aligned, General
# format" style already used throughout rows 1-3 (styles.xml index 2).
$rightCols = @("A","B","C","D","E","F","G","H","I","K","L","M","N","O","P","Q",
               "S","T","U","V","W","X","Y","Z",
               "AA","AB","AC","AD","AE","AF","AG","AH","AI")
foreach ($col in $rightCols) {
    $ws.Range($col + "4").HorizontalAlignment = $xlRight
}
# Column R keeps the "left aligned" style used in rows 1-3 (index 1).
$ws.Range("R4").HorizontalAlignment = $xlLeft

# --- values ---------------------------------------------------------------
# Copy/PasteSpecial (values only) preserves the source's underlying storage
# type, so text that looks numeric (e.g. the leading-zero phone number in
# H) is carried over as text instead of being reinterpreted as a number.
$srcCols = @("B","C","D","E","F","G","H","J","K","L","M","N","O")
foreach ($col in $srcCols) {
    $ws.Range($col + "3").Copy()
    $ws.Range($col + "4").PasteSpecial($xlPasteValues)
}

# A4 is the new row's own sequence number (3), not copied from A3 (2).
$ws.Range("A4").Value = 3

# --- move the saved selection to J2, matching the sheetView in the diff --
[void]$ws.Range("J2").Select()
